$wb = $excel.ActiveWorkbook

# Helper: set a cell's value as genuine text (not auto-converted to a number)
# while preserving the cell's existing style index. We do this by writing a
# text formula (="...") and then converting it in-place to a static value via
# Copy + PasteSpecial(xlPasteValues). This keeps the pre-existing style
# (borders/bold/alignment) untouched, unlike setting .NumberFormat which would
# allocate a brand-new style entry.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy($range)
    $range.PasteSpecial(-4163) # xlPasteValues
}

# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" -- fix mislabeled E1 header (2050) and drop the
# trailing "Total" row (row 13).
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Rows.Item(13).Delete()
    Set-TextValue $ws.Range("E1") "2050"
}

# Sheet 4: "Potencia Incremental - SIN(MW)" -- same fix, but the header uses a
# year-range label ("2041-2050") instead of a plain year.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(13).Delete()
Set-TextValue $ws4.Range("E1") "2041-2050"

# Sheet 5: "Emissoes Totais (MtCO2eq)" -- only the mislabeled E1 header needs
# fixing; this sheet never had a "Total" row.
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("E1") "2050"

# Sheet 6: "Custo Total (bilhões de R$)" -- no E1 header on this sheet; just
# drop the trailing "Total" row (row 4).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
